$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$meta.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicate "Contact" row (row 10), shifting everything below up by one
$meta.Rows.Item(10).Delete()

# Publisher value (row 9) now gets a real publisher name
$meta.Range("B9").Value = "Alvearie Team"

# The row that used to be the (duplicate) Contact row is now Jurisdiction
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root extension row: Short / Definition now reflect the real title/description
$elements.Range("K2").Value = "Child Organization Hierarchy Level Description"
$elements.Range("L2").Value = "Description of the level of the child practitioner within the organinzational hierarchy"

# Column K (Short) grows wider to fit the new text (bestFit recalculation)
$elements.Columns.Item(11).ColumnWidth = 43.2
